# Update Argent (silver) price history sheet:
#  - B2 ("03/02/2025" row) price revised from 31,53 to 31,577
#  - New row 4 added for 04/02/2025 with the previous price of 31,53
#
# The source values are plain text (French-style decimal commas and
# dd/mm/yyyy dates) rather than real numbers/dates, so each target cell is
# explicitly formatted as Text before the value is written - this stops
# Excel's automatic number/date recognition from turning "31,577" into the
# number 31577 or "04/02/2025" into a date serial.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Revise the existing 03/02/2025 price.
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "31,577"

# Append the new 04/02/2025 row.
$ws.Range("A4").NumberFormat = "@"
$ws.Range("A4").Value = "04/02/2025"

$ws.Range("B4").NumberFormat = "@"
$ws.Range("B4").Value = "31,53"
